# Update computed values in Sheet1 for rows 2-25, columns B, C, D, F, G, L, O
# per commit 'case with 380 kV done'
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cols = @('B','C','D','F','G','L','O')

$data = @(
    @(1.406305914058237, 0.3080582067633486, 0.02522507647812944, 0.692197308576894, 0.002420429848090936, 0.3014673189856438, 2.355552748963788),
    @(1.267486375943577, 0.2955657763784814, 0.0237497511240008, 0.6916207774384944, 0.002423550382351923, 0.2901248941801242, 2.369412810402082),
    @(1.18227271276038, 0.2879004649577155, 0.02284106791326224, 0.6918918423909446, 0.002425568306331232, 0.2833379650290766, 2.38006358789346),
    @(1.147554754435419, 0.2847783730288711, 0.02247008952500806, 0.6921590727995479, 0.002426416329588523, 0.2806168046047048, 2.384940698703673),
    @(1.141790356792797, 0.2842600567945794, 0.02240844846156165, 0.6922129027559052, 0.002426558697915185, 0.2801676489233955, 2.385782920073552),
    @(1.181804462338391, 0.2878583525116483, 0.0228360674945165, 0.691894812142948, 0.002425579639012169, 0.2833010860924929, 2.380127190634568),
    @(1.35843734846668, 0.3037499374159722, 0.02471698759708829, 0.6918685431868994, 0.002421484705534844, 0.2975196149965598, 2.359886532134112),
    @(1.70493074619236, 0.3349420315788336, 0.02838197612705073, 0.6967961153916917, 0.002414259522464376, 0.3268127416539244, 2.33724726552731),
    @(1.959518968247835, 0.3578632555122567, 0.03105918714922495, 0.703480557918482, 0.002409436876486115, 0.3492016812680987, 2.331106556592147),
    @(2.075333049838832, 0.3682891068921208, 0.03227354665310855, 0.7071931086704524, 0.002407347322000854, 0.3595770672998526, 2.330612244794082),
    @(2.119187620851108, 0.3722366847980538, 0.0327328638325497, 0.7086960242149019, 0.002406570977775326, 0.3635334418081868, 2.330757309883751),
    @(2.109742858651657, 0.3713865287621729, 0.03263396587873046, 0.7083680207820748, 0.002406737514783114, 0.3626801446632157, 2.330711264554083),
    @(2.078941036231811, 0.3686138875579843, 0.03231134582638617, 0.7073148067516541, 0.002407283152835597, 0.3599020104159791, 2.330617511226762),
    @(2.060073755236147, 0.3669154951828659, 0.03211366129069404, 0.7066823353826663, 0.00240761931367547, 0.3582038978288011, 2.330603401951294),
    @(1.951950103920922, 0.3571818528491519, 0.03097975266635444, 0.7032514904835807, 0.002409575526326334, 0.3485274631135553, 2.331185265514335),
    @(1.885618618251101, 0.351210072123564, 0.03028321508832477, 0.7013191570138844, 0.00241080225937435, 0.3426401157529995, 2.332132248877087),
    @(1.847466666536661, 0.3477751734332628, 0.02988225572165959, 0.7002709325026046, 0.002411517664471684, 0.3392718017769454, 2.332893235719609),
    @(1.834549154660294, 0.3466121699047449, 0.0297464420792366, 0.6999268629325428, 0.002411761577383157, 0.3381344272358433, 2.333187993753342),
    @(1.892679714223164, 0.3518457897349379, 0.03035739700015938, 0.7015183126811451, 0.002410670655700783, 0.3432649776916037, 2.332009042020871),
    @(2.087988338685079, 0.3694282945250222, 0.03240612188593417, 0.7076215236367602, 0.002407122481466307, 0.3607172707161368, 2.33063601863779),
    @(2.215623042114771, 0.3809166255877301, 0.03374195411602443, 0.7121762343891049, 0.002404890504305856, 0.3722833198608271, 2.331676034634938),
    @(2.147503571987556, 0.3747854471322398, 0.03302929086459017, 0.7096933691474305, 0.002406073819331144, 0.3660956510342004, 2.330943149096555),
    @(1.889487445859572, 0.3515583868961301, 0.03032386094601947, 0.701428079130423, 0.002410730122138383, 0.3429824265300709, 2.332064069401355),
    @(1.611188037029365, 0.3265020420094231, 0.02739313356759254, 0.6949270506254379, 0.002416128470682645, 0.3187364039770273, 2.341536183995032),
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    for ($j = 0; $j -lt $cols.Count; $j++) {
    $ws.Range("$($cols[$j])$row").Value = $data[$i][$j]
    }
}
